# Insert a new data row for 2026/02/14 08:00 (Saturday) at row 807 of Sheet1,
# pushing the existing rows 807-848 down to 808-849.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 807; everything below shifts down by one.
$ws.Rows.Item(807).Insert()

# Column A holds the date as plain text (e.g. "2026/12/29"), not a real Excel
# date, matching how the rest of the sheet stores it. A leading apostrophe
# forces the literal text instead of Excel's auto date-detection; resetting
# the style back to "Normal" afterwards removes the quote-prefix styling so
# the cell matches its neighbours exactly.
$ws.Cells.Item(807, 1).Value = "'2026/02/14"
$ws.Cells.Item(807, 1).Style = "Normal"

$ws.Cells.Item(807, 2).Value = "土"
$ws.Cells.Item(807, 3).Value = 8
$ws.Cells.Item(807, 4).Value = 201
